$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.900.32"
$ws.Range("E2").Value = "  +5.86%  "

# Row 3
$ws.Range("D3").Value = "2.251.51"
$ws.Range("E3").Value = "  +4.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.92"
$ws.Range("E5").Value = "  +0.93%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.43"
$ws.Range("E7").Value = "  -2.98%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.405"
$ws.Range("E9").Value = "  +3.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.61"
$ws.Range("E10").Value = "  +0.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0881"
$ws.Range("E11").Value = "  +3.87%  "

# Row 12
$ws.Range("E12").Value = "  +0.75%  "

# Row 13
$ws.Range("D13").Value = "2.585.07"
$ws.Range("E13").Value = "  +4.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.93"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.77"
$ws.Range("E15").Value = "  -0.89%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.804"
$ws.Range("E16").Value = "  -0.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.63"
$ws.Range("E17").Value = "  +2.64%  "

# Row 18
$ws.Range("D18").Value = "2.248.06"
$ws.Range("E18").Value = "  +4.13%  "

# Row 19
$ws.Range("D19").Value = "41.760.12"
$ws.Range("E19").Value = "  +5.55%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.31"
$ws.Range("E20").Value = "  +1.99%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +6.65%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  +0.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.09"
$ws.Range("E23").Value = "  +8.79%  "

# Row 24
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.39"
$ws.Range("E25").Value = "  +0.78%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  +0.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("E27").Value = "  +0.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("E28").Value = "  +3.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.06"
$ws.Range("E29").Value = "  -2.28%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.20"
$ws.Range("E30").Value = "  +2.27%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  +2.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  +4.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("E34").Value = "  +9.53%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  +1.45%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0628"
$ws.Range("E36").Value = "  +1.60%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("E37").Value = "  +5.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.75"
$ws.Range("E38").Value = "  -3.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("E39").Value = "  -0.04%  "

# Row 40
$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000250"
$ws.Range("E40").Value = "  +30.40%  "

# Row 41
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.29%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.93"
$ws.Range("E42").Value = "  +4.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0239"
$ws.Range("E43").Value = "  +5.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.76"
$ws.Range("E44").Value = "  +13.25%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.73"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0979"
$ws.Range("E46").Value = "  +5.80%  "

# Row 47
$ws.Range("D47").Value = "1.488.20"
$ws.Range("E47").Value = "  -1.68%  "

# Row 48
$ws.Range("E48").Value = "  -1.10%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.65"
$ws.Range("E49").Value = "  -5.64%  "

# Row 50
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  +0.56%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("E51").Value = "  -0.16%  "
